# Change INGESTION_RULE ('NONE') to 'STRING' for the rows whose rule was
# previously 'NONE' in the mapping-sheet worksheet, and keep the
# pre-computed O column (paste-as-values snapshot of column N) in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping-sheet")

# Rows (2-41) whose INGESTION_RULE (column L) value is currently 'NONE'.
$rows = @(2, 3, 6, 7, 9, 10, 11, 16, 17, 18, 19, 24, 26, 27, 29, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41)

foreach ($r in $rows) {
    $lCell = $ws.Cells.Item($r, 12)   # column L = INGESTION_RULE
    if ($lCell.Value() -eq "NONE") {
        $lCell.Value = "STRING"
    }
    # Column N recalculates automatically (CONCATENATE formula). Column O
    # is a static snapshot of that same text, so refresh it to match.
    $nCell = $ws.Cells.Item($r, 14)   # column N
    $oCell = $ws.Cells.Item($r, 15)   # column O
    $oCell.Value = $nCell.Value2
}

# Update the hidden _xlnm._FilterDatabase defined name to cover the full
# table range A1:O68 instead of just K1:K68.
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "='mapping-sheet'!`$A`$1:`$O`$68"

# Refresh the view position/selection to match the post-edit state.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("O18").Select() | Out-Null
